$wb = $excel.ActiveWorkbook

# Add the new "Storedata" worksheet right after "userdata"
$afterSheet = $wb.Worksheets.Item("userdata")
$ws = $wb.Worksheets.Add($null, $afterSheet)
$ws.Name = "Storedata"

# Header row
$ws.Range("A1").Value = "OrderId"
$ws.Range("B1").Value = "petId"
$ws.Range("C1").Value = "Quantity"
$ws.Range("D1").Value = "Shipdate"
$ws.Range("E1").Value = "status"
$ws.Range("F1").Value = "complete"

# Header styling: solid fill, Gold Accent4 Lighter 60%
$headerRange = $ws.Range("A1:F1")
$headerRange.Interior.ThemeColor = 8

# Data rows
$ws.Range("A2").Value = 123
$ws.Range("B2").Value = 8
$ws.Range("C2").Value = 43
$ws.Range("D2").Value = "2023-11-03T05:42:06.082Z"
$ws.Range("E2").Value = "placed"
$ws.Range("F2").Value = $true

$ws.Range("A3").Value = 15
$ws.Range("B3").Value = 4
$ws.Range("C3").Value = 89
$ws.Range("D3").Value = "2023-11-03T05:42:06.082Z"
$ws.Range("E3").Value = "pending"
$ws.Range("F3").Value = $false

$ws.Range("A4").Value = 75
$ws.Range("B4").Value = 6
$ws.Range("C4").Value = 7
$ws.Range("D4").Value = "2023-11-03T05:42:06.082Z"
$ws.Range("E4").Value = "conformed"
$ws.Range("F4").Value = $true

$ws.Range("A5").Value = 18
$ws.Range("B5").Value = 7
$ws.Range("C5").Value = 1
$ws.Range("D5").Value = "2023-11-03T05:42:06.082Z"
$ws.Range("E5").Value = "placed"
$ws.Range("F5").Value = $true

$ws.Range("A6").Value = 12
$ws.Range("B6").Value = 3
$ws.Range("C6").Value = 3
$ws.Range("D6").Value = "2023-11-03T05:42:06.082Z"
$ws.Range("E6").Value = "cancelled"
$ws.Range("F6").Value = $false

# Column widths for D (Shipdate) and E (status)
$ws.Columns.Item(4).ColumnWidth = 32
$ws.Columns.Item(5).ColumnWidth = 14.15

# Make Storedata the active sheet/tab, with the same selection the author left it in
$ws.Activate()
$ws.Range("D10").Select()
